# microgrid.py not calculating m2 loads as expected
#
# Duplicate the "12222020_Lehigh" sheet into a brand-new first sheet named
# "12282020_Lehigh" that carries corrected m2 (Min/Avg/Max) load figures,
# while leaving the original "12222020_Lehigh" sheet (now second in tab
# order) untouched aside from its view/selection.

$wb = $excel.ActiveWorkbook

# 1. Duplicate the first worksheet ("12222020_Lehigh") placing the copy
#    before it, so it becomes the new first sheet/tab.
$src = $wb.Worksheets.Item(1)
$src.Copy($src)
$newSheet = $wb.Worksheets.Item(1)
$newSheet.Name = "12282020_Lehigh"

# 2. The new sheet's column headers no longer reference the 12/22/2020
#    snapshot date - make them generic.
$newSheet.Range("I10").Value = "Average Load (kW) "
$newSheet.Range("J10").Value = "Max Load (kW) "
$newSheet.Range("K10").Value = "Min Load (kW) "

# 3. Corrected m2 load numbers (row 15, the "m2" microgrid).
$newSheet.Range("I15").Value = 163
$newSheet.Range("J15").Value = 389
$newSheet.Range("K15").Value = 53

# 4. Corrected average outage-survival figure for the full grid comparison.
$newSheet.Range("Y28").Value = 44

# 5. Freeze panes on the new sheet (header row + first column) and leave
#    the selection on D6, matching the reviewed layout.
$newSheet.Activate()
$newSheet.Range("B11").Select()
$excel.ActiveWindow.FreezePanes = $true
$newSheet.Range("D6").Select()

# 6. The old "12222020_Lehigh" sheet (now 2nd tab) is no longer the active
#    tab; just scroll its view over to column H.
$oldSheet = $wb.Worksheets.Item(2)
$oldSheet.Activate()
$oldSheet.Range("H1").Select()

# 7. Make the new sheet the active/selected tab again.
$newSheet.Activate()
